$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8 (shifts old rows 8-9 down to 9-10)
$ws.Rows.Item(8).Insert()

# Fill the new row 8 with the "Possible_Problem" branch data
# (same Node1 text as the rows below it, new Relationship/Node2 values)
$ws.Cells.Item(8, 1).Value = "Problem:Does driving the vehicle alleviate the problem? (Please answer as: Yes, No)"
$ws.Cells.Item(8, 2).Value = "Possible_Problem"
$ws.Cells.Item(8, 3).Value = "Possible_Problem:30% Restricted Heater Core`n20%Thermostat`n15% HVAC Door Actuators`n10% Restricted Radiator Coolant Flow`n10% Low Coolant Level`n5% HVAC Control Unit`n5% Water Pump`n5% Cylinder Head Gasket"

# Apply wrap text style to the new C8 cell (matches style used by other Node2 cells)
$ws.Cells.Item(8, 3).WrapText = $true

# Restore the explicit row height for the new row 8
$ws.Rows.Item(8).RowHeight = 409.6

# Update selection / active cell to match the saved view state
$ws.Range("C8").Select()
